$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.056.19'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.892.77'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.69'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07262'
$ws.Range("E9").Value = '  +0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.09'
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8991'
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08167'
$ws.Range("E12").Value = '  +6.34%  '
$ws.Range("D13").Value = '1.942.23'
$ws.Range("E13").Value = '  +1.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '96.21'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008578'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("E18").Value = '  +1.08%  '
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '27.107.26'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.077'
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.68'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.405'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '148.62'
$ws.Range("E24").Value = '  +2.05%  '
$ws.Range("E25").Value = '  -0.74%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.733'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '114.94'
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.780'
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.840'
$ws.Range("E30").Value = '  -2.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09223'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05035'
$ws.Range("E32").Value = '  -0.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7893'
$ws.Range("E33").Value = '  -2.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.209'
$ws.Range("E34").Value = '  -2.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.429'
$ws.Range("E35").Value = '  +3.42%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.978'
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.603'
$ws.Range("E37").Value = '  +0.64%  '
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01985'
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.044'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.551'
$ws.Range("E42").Value = '  -0.85%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.25'
$ws.Range("E43").Value = '  -2.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1514'
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4867'
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.04'
$ws.Range("E47").Value = '  -1.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.623'
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.14'
$ws.Range("E49").Value = '  +1.71%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.49'
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05930'
$ws.Range("E51").Value = '  -0.02%  '
